$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix team name mismatches
$ws.Range("A18").Value = "North Melbourne"
$ws.Range("A19").Value = "North Melbourne"
$ws.Range("A14").Value = "Gold Coast"

# Update the selected cell to match the saved state
$ws.Range("A14").Select()
